$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$helper = $ws.Range("Z1")

# Row 2 - Bitcoin
$ws.Range("D2").Value = "67.812.05"

$ws.Range("E2").Value = "  -2.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.559.40"

$ws.Range("E3").Value = "  -3.44%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$helper.NumberFormat = "@"
$helper.Value = "616.18"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E5").Value = "  -7.32%  "

# Row 6 - Solana
$helper.NumberFormat = "@"
$helper.Value = "153.76"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E6").Value = "  -3.95%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.556.47"

$ws.Range("E7").Value = "  -3.43%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -2.16%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.17%  "

# Row 11 - Toncoin
$helper.NumberFormat = "@"
$helper.Value = "6.91"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E11").Value = "  -3.27%  "

# Row 12 - Cardano
$helper.NumberFormat = "@"
$helper.Value = "0.432"
$helper.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E12").Value = "  -1.94%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -3.47%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.159.06"

$ws.Range("E14").Value = "  -3.51%  "

# Row 15 - Avalanche
$helper.NumberFormat = "@"
$helper.Value = "32.11"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E15").Value = "  -2.05%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.559.68"

$ws.Range("E16").Value = "  -3.02%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "67.817.24"

$ws.Range("E17").Value = "  -2.54%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.97%  "

# Row 19 - Polkadot
$helper.NumberFormat = "@"
$helper.Value = "6.42"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E19").Value = "  -0.61%  "

# Row 20 - Chainlink
$helper.NumberFormat = "@"
$helper.Value = "15.69"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E20").Value = "  -2.52%  "

# Row 21 - BitcoinCash
$helper.NumberFormat = "@"
$helper.Value = "453.93"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E21").Value = "  -3.08%  "

# Row 22 - Uniswap
$helper.NumberFormat = "@"
$helper.Value = "9.61"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E22").Value = "  -1.48%  "

# Row 23 - Polygon
$helper.NumberFormat = "@"
$helper.Value = "0.646"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E23").Value = "  +0.14%  "

# Row 24 - Litecoin
$helper.NumberFormat = "@"
$helper.Value = "77.62"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E24").Value = "  -2.80%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.702.92"

$ws.Range("E25").Value = "  -3.39%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.09%  "

# Row 27 - InternetComputer(DFINITY)
$helper.NumberFormat = "@"
$helper.Value = "10.65"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E27").Value = "  -2.52%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -7.60%  "

# Row 29 - RenderToken
$helper.NumberFormat = "@"
$helper.Value = "8.40"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E29").Value = "  -6.89%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -3.95%  "

# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -3.73%  "

# Row 32 - Binance-PegBSC-USD
$ws.Range("E32").Value = "  +0.07%  "

# Row 33 - EthereumClassic
$helper.NumberFormat = "@"
$helper.Value = "25.96"
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E33").Value = "  -2.84%  "

# Row 34 - ImmutableX
$ws.Range("E34").Value = "  -4.48%  "

# Row 35 - NEARProtocol
$helper.NumberFormat = "@"
$helper.Value = "6.22"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E35").Value = "  -3.81%  "

# Row 36 - Kaspa (was RenzoRestakedETH)
$ws.Range("B36").Value = "Kaspa"

$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"

$helper.NumberFormat = "@"
$helper.Value = "0.157"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E36").Value = "  -4.18%  "

# Row 37 - RenzoRestakedETH (was Kaspa)
$ws.Range("B37").Value = "RenzoRestakedETH"

$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"

$ws.Range("D37").Value = "3.558.49"

$ws.Range("E37").Value = "  -3.26%  "

# Row 38 - Aptos
$helper.NumberFormat = "@"
$helper.Value = "8.09"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E38").Value = "  -4.03%  "

# Row 40 - FirstDigitalUSD
$helper.NumberFormat = "@"
$helper.Value = "0.999"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E40").Value = "  -0.08%  "

# Row 41 - Monero
$helper.NumberFormat = "@"
$helper.Value = "176.74"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E41").Value = "  -1.16%  "

# Row 42 - Hedera
$helper.NumberFormat = "@"
$helper.Value = "0.0886"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E42").Value = "  -1.55%  "

# Row 43 - Filecoin
$helper.NumberFormat = "@"
$helper.Value = "5.62"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E43").Value = "  -7.47%  "

# Row 44 - Stacks
$helper.NumberFormat = "@"
$helper.Value = "2.10"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E44").Value = "  -6.52%  "

# Row 45 - Mantle
$helper.NumberFormat = "@"
$helper.Value = "0.895"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E45").Value = "  -4.07%  "

# Row 46 - InjectiveProtocol (was OKB)
$ws.Range("B46").Value = "InjectiveProtocol"

$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"

$helper.NumberFormat = "@"
$helper.Value = "29.07"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E46").Value = "  +6.25%  "

# Row 47 - OKB (was InjectiveProtocol)
$ws.Range("B47").Value = "OKB"

$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

$helper.NumberFormat = "@"
$helper.Value = "46.25"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E47").Value = "  -1.56%  "

# Row 48 - dogwifhat
$helper.NumberFormat = "@"
$helper.Value = "2.60"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Clear()

$ws.Range("E48").Value = "  -5.35%  "

# Row 49 - Cosmos
$ws.Range("E49").Value = "  -1.71%  "

# Row 50 - ONDO
$ws.Range("E50").Value = "  -5.96%  "

# Row 51 - SuiNetwork
$ws.Range("E51").Value = "  -4.66%  "
